# Updates cryptos list price/volume figures (and corrects the
# Monero/Fetch.AI row ordering) as produced by the scheduled
# GitHub Actions refresh job.
# Numeric-looking price strings are prefixed with a leading apostrophe
# so Excel stores them as text (matching the workbook's original
# inlineStr/text cells) instead of silently converting them to
# floating point numbers and losing trailing zeros / exact formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.866.68'
$ws.Range("E2").Value = '  -1.50%  '
$ws.Range("D3").Value = '3.075.06'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''520.04'
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").Value = '''136.05'
$ws.Range("E6").Value = '  -3.35%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").Value = '3.073.35'
$ws.Range("D9").Value = '''0.450'
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("D10").Value = '''7.34'
$ws.Range("E10").Value = '  +2.89%  '
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").Value = '3.601.39'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '''25.23'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("E16").Value = '  -2.10%  '
$ws.Range("D17").Value = '56.974.08'
$ws.Range("E17").Value = '  -1.40%  '
$ws.Range("D18").Value = '3.072.01'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '''5.87'
$ws.Range("E19").Value = '  -3.44%  '
$ws.Range("D20").Value = '''12.44'
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("D21").Value = '''7.83'
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").Value = '''346.14'
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").Value = '''68.28'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").Value = '''0.497'
$ws.Range("E26").Value = '  -2.63%  '
$ws.Range("D27").Value = '''0.167'
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '0.0₃0861'
$ws.Range("E29").Value = '  -5.79%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("D31").Value = '''7.26'
$ws.Range("E31").Value = '  +0.67%  '
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").Value = '''5.86'
$ws.Range("E33").Value = '  -7.86%  '
$ws.Range("D34").Value = '''20.72'
$ws.Range("E34").Value = '  -0.86%  '
$ws.Range("D35").Value = '''4.91'
$ws.Range("E35").Value = '  +6.47%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = '''1.14'
$ws.Range("E36").Value = '  -3.74%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '''158.84'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").Value = '''5.98'
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("D39").Value = '''25.65'
$ws.Range("E39").Value = '  -1.15%  '
$ws.Range("E40").Value = '  -1.91%  '
$ws.Range("D41").Value = '''0.0651'
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").Value = '''4.00'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").Value = '''0.688'
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("D45").Value = '2.377.38'
$ws.Range("E45").Value = '  +4.39%  '
$ws.Range("D46").Value = '''36.55'
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = '3.112.05'
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").Value = '''0.0263'
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").Value = '''0.954'
$ws.Range("E50").Value = '  -3.82%  '
$ws.Range("D51").Value = '''5.93'
$ws.Range("E51").Value = '  -2.25%  '
